# Update the shared "Заповедь/Посыл" text (stored as an encrypted blob) that
# is repeated in cells B3, B5 and B7. Setting the same new value into every
# cell that previously pointed at the old shared string lets the writer
# collapse them back into a single shared-string table entry, matching how
# the workbook was re-saved upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "U2FsdGVkX1+OSmLAeytfRFMzwftOg7BmHQT6GYLtwb17ew7I/rZrHXM2KOujqXwc/aBvndfAz6U5ltHv8g67vD7DOQSPPIUTYK2XM/4+XGMNqlxtciFXSfm7nytpBIQOWC3qB+e8ZqUaPUkacfkVQa3MwERnX8LBaNdq3C5bGlqDxO951iCgbOLjn+LuohvCFk0+u4F9shF8B/yPl2YJOZpzSQaWk8hmd7zfHrNOn6aUBWGJi7EJicqJX19/EsyqniFCxBAPKnUs6R3AanXuIXsKGPAwomOLZtwc2yQm6+GHPRacQ0HoT3LEfWC6VjjRWvNWV31EO26KUM3dWj1VRUULSQjlUzOWgYr+8solE0oMz8sxIlda08Bc0ERV4w3t"

$ws.Range("B3").Value = $newValue
$ws.Range("B5").Value = $newValue
$ws.Range("B7").Value = $newValue

# Move the active selection from B11 to B10, as in the saved workbook.
$ws.Range("B10").Select()
